$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239 - this shifts the existing rows
# 239..340 down to 240..341, preserving their original values.
$ws.Rows.Item(239).Insert()

# Populate the newly-inserted row 239 with the new record's data.
$ws.Range("A239").Value = 3
$ws.Range("B239").Value = "Femacal de La Calera"
$ws.Range("C239").Value = "Coquimbo"
$ws.Range("D239").Value = 44992
$ws.Range("E239").Value = 5
$ws.Range("F239").Value = "Fruta"
$ws.Range("G239").Value = 100101
$ws.Range("H239").Value = "Berries"
$ws.Range("I239").Value = 100101001
$ws.Range("J239").Value = "Arándano (blue)"
$ws.Range("K239").Value = "Sin especificar"
$ws.Range("L239").Value = "Primera"
$ws.Range("M239").Value = 30
$ws.Range("N239").Value = 3600
$ws.Range("O239").Value = 3600
$ws.Range("P239").Value = 3600
$ws.Range("Q239").Value = "`$/bandeja 2 kilos"
$ws.Range("R239").Value = "Provincia de Curicó"
$ws.Range("S239").Value = 1800
$ws.Range("T239").Value = 2
